{"js": "// Beautify html CV with strapdown.js.\n//\n// The top contact-info list (numId=2) is reordered/pruned:\n//   \u6027\u522b: \u7537                               -> becomes \"\u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com\"\n//   \u51fa\u751f\u65e5\u671f: 1985 \u5e74 11 \u6708                  -> unchanged, but a new \"\u79fb\u52a8\u7535\u8bdd\" item lands right before it\n//   \u5c45\u4f4f\u5730: \u4e0a\u6d77\u5e02                          -> removed\n//   \u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com        -> removed (its text now lives in the first bullet)\n//   \u79fb\u52a8\u7535\u8bdd: 13764131714                   -> removed (its text now lives in the new bullet before \u51fa\u751f\u65e5\u671f)\n//\n// Also drop the trailing stray \"*\" paragraph at the very end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst texts = paragraphs.items.map((p) => p.text.trim());\n\nfunction findIndex(match) {\n  const i = texts.indexOf(match);\n  if (i === -1) {\n    throw new Error(`Paragraph with text ${JSON.stringify(match)} not found`);\n  }\n  return i;\n}\n\nconst idxGender = findIndex(\"\u6027\u522b: \u7537\");\nconst idxBirth = findIndex(\"\u51fa\u751f\u65e5\u671f: 1985 \u5e74 11 \u6708\");\nconst idxCity = findIndex(\"\u5c45\u4f4f\u5730: \u4e0a\u6d77\u5e02\");\nconst idxEmail = findIndex(\"\u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com\");\nconst idxPhone = findIndex(\"\u79fb\u52a8\u7535\u8bdd: 13764131714\");\n\n// 1. \"\u6027\u522b: \u7537\" -> \"\u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com\"\nparagraphs.items[idxGender].insertText(\"\u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com\", \"Replace\");\n\n// 2. Insert a new \"\u79fb\u52a8\u7535\u8bdd: 13764131714\" bullet right before \"\u51fa\u751f\u65e5\u671f...\" \u2014\n//    it inherits the numbered-list paragraph formatting from that paragraph.\nparagraphs.items[idxBirth].insertParagraph(\"\u79fb\u52a8\u7535\u8bdd: 13764131714\", \"Before\");\n\n// 3. Drop the now-duplicated \"\u5c45\u4f4f\u5730\", \"\u7535\u5b50\u90ae\u4ef6\" and \"\u79fb\u52a8\u7535\u8bdd\" paragraphs.\nparagraphs.items[idxCity].delete();\nparagraphs.items[idxEmail].delete();\nparagraphs.items[idxPhone].delete();\n\n// 4. Remove the trailing \"*\" paragraph at the end of the document.\nconst lastParagraphs = body.paragraphs;\nlastParagraphs.load(\"text\");\nawait context.sync();\n\nconst items = lastParagraphs.items;\nconst last = items[items.length - 1];\nif (last.text.trim() === \"*\") {\n  last.delete();\n}\n\nawait context.sync();\n", "ps1": "# Beautify html CV with strapdown.js.\n#\n# The top contact-info list (numId=2) is reordered/pruned:\n#   \u6027\u522b: \u7537                               -> becomes \"\u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com\"\n#   \u51fa\u751f\u65e5\u671f: 1985 \u5e74 11 \u6708                  -> unchanged, but a new \"\u79fb\u52a8\u7535\u8bdd\" item lands right before it\n#   \u5c45\u4f4f\u5730: \u4e0a\u6d77\u5e02                          -> removed\n#   \u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com        -> removed (its text now lives in the first bullet)\n#   \u79fb\u52a8\u7535\u8bdd: 13764131714                   -> removed (its text now lives in the new bullet before \u51fa\u751f\u65e5\u671f)\n#\n# Also drop the trailing stray \"*\" paragraph at the very end of the document.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $text, $startAt) {\n    $count = $doc.Paragraphs.Count\n    for ($i = $startAt; $i -le $count; $i++) {\n        $t = $doc.Paragraphs($i).Range.Text.TrimEnd()\n        if ($t -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\n$idxGender = Find-ParagraphIndex $d \"\u6027\u522b: \u7537\" 1\n$idxBirth  = Find-ParagraphIndex $d \"\u51fa\u751f\u65e5\u671f: 1985 \u5e74 11 \u6708\" 1\n\n# 1. \"\u6027\u522b: \u7537\" -> \"\u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com\"\n$d.Paragraphs($idxGender).Range.Text = \"\u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com\"\n\n# 2. Insert a new \"\u79fb\u52a8\u7535\u8bdd: 13764131714\" bullet right before \"\u51fa\u751f\u65e5\u671f...\" \u2014\n#    it inherits the numbered-list paragraph formatting from that paragraph.\n$d.Paragraphs($idxBirth).Range.InsertParagraphBefore()\n$d.Paragraphs($idxBirth).Range.Text = \"\u79fb\u52a8\u7535\u8bdd: 13764131714\"\n\n# 3. Drop the now-duplicated \"\u5c45\u4f4f\u5730\", \"\u7535\u5b50\u90ae\u4ef6\" and \"\u79fb\u52a8\u7535\u8bdd\" paragraphs.\n#    They all sit after the freshly-rewritten bullets above (idxBirth now\n#    holds the new \"\u79fb\u52a8\u7535\u8bdd\" paragraph, and the original \"\u51fa\u751f\u65e5\u671f\" got\n#    pushed one slot further down), so search past that point to land on\n#    the stale duplicates rather than the paragraphs we just wrote.\n$searchStart = $idxBirth + 2\n\n$idxCity = Find-ParagraphIndex $d \"\u5c45\u4f4f\u5730: \u4e0a\u6d77\u5e02\" $searchStart\n$d.Paragraphs($idxCity).Range.Delete()\n\n$idxEmail = Find-ParagraphIndex $d \"\u7535\u5b50\u90ae\u4ef6: logan.zhou.cn@gmail.com\" $searchStart\n$d.Paragraphs($idxEmail).Range.Delete()\n\n$idxPhone = Find-ParagraphIndex $d \"\u79fb\u52a8\u7535\u8bdd: 13764131714\" $searchStart\n$d.Paragraphs($idxPhone).Range.Delete()\n\n# 4. Remove the trailing \"*\" paragraph at the end of the document.\n$last = $d.Paragraphs.Count\n$lastText = $d.Paragraphs($last).Range.Text.TrimEnd()\nif ($lastText -eq \"*\") {\n    $d.Paragraphs($last).Range.Delete()\n}\n"}
